$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Producer name change
$ws.Range("B1").Value = "La ferme de Maurice"

# --- Product block 1 (rows 7-10): was "Pain complet" / A la pièce -> "Tomates grappe" / Au poids (Kg)
# (A7:A8 is a merged pair; writing the visible anchor cell A7 is what Excel allows/shows.)
$ws.Range("A7").Value = "Tomates grappe"
$ws.Range("B7").Value = "Au poids (Kg)"
$ws.Range("C7").Value = 2
$ws.Range("B9").Value = 1

# --- Product block 2 (rows 12-15): was "Tomates grappe" / Au poids (Kg) -> "Pomme de terre" / Au poids (Kg)
$ws.Range("A12").Value = "Pomme de terre"
$ws.Range("B12").Value = "Au poids (Kg)"
$ws.Range("C12").Value = 2
$ws.Range("B14").Value = 1

# --- Product block 3 (rows 17-20): was "Pomme de terre" / Au poids (Kg) -> "Radis" / A la pièce
$ws.Range("A17").Value = "Radis"
$ws.Range("B17").Value = "A la pièce"
$ws.Range("C17").Value = 4
$ws.Range("B19").Value = 1

# --- Product block 4 (rows 22-25): was "Salade" / A la pièce -> stays "Salade" / A la pièce
$ws.Range("A22").Value = "Salade"
$ws.Range("B22").Value = "A la pièce"
$ws.Range("C22").Value = 1
$ws.Range("B24").Value = 1

# --- Remove product block 5 (rows 27-30, "Radis" / A la pièce) and the blank separator row 31.
# This shifts the old grand-total row (32) up to row 27 automatically, adjusting the
# mergeCells list and the SUM formula's row references.
$ws.Range("A27:D31").EntireRow.Delete()

# The deletion leaves a #REF! term where the deleted C30 used to be; restore the
# intended formula that now only sums the 4 remaining subtotal rows.
$ws.Range("C27").Formula = "=C10+C15+C20+C25"

Write-Output "edit applied"
